$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the 47177279 row (row 4)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-01 10:52:31"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime" for the 47177279 row (row 4)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-01 10:52:26"
$wsZhCn.Range("K4").Value = "2016-09-01 10:53:05"

# de-de sheet: "Correspond Handback DateTime" for the 47177279 row (row 4)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-09-01 10:53:17"
